$d = $word.ActiveDocument

# --- 1) First paragraph: add trailing spaces + red "(This is a change ..." runs ---
$d.Content.Find.Execute(
    "This is a Microsoft word document.", $true, $false, $false, $false, $false,
    $true, 1, $false, "This is a Microsoft word document.  ", 2) | Out-Null

$enDash = [char]0x2013
$redColor = 192  # RGB(0xC0,0x00,0x00) = C00000

$insertPos = $d.Paragraphs(1).Range.End - 1

$r1 = $d.Range($insertPos, $insertPos)
$r1.InsertAfter("(This is a change " + $enDash + " Ve")
$r1.Font.Color = $redColor

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("rsion for branch alternate")
$r2.Font.Color = $redColor

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter(")")
$r3.Font.Color = $redColor

# --- 2) Append a new empty, shaded paragraph at the very end of the body ---
$endOfDoc = $d.Content.End
$tail = $d.Range($endOfDoc, $endOfDoc)
$shadedParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
[void]$tail.InsertXML($shadedParaXml)
